# =====================================================================
# Adds a "2022-Q3" sheet (placed right after "总计") with fund-holding
# data, and updates the "总计" summary sheet with a new first data row
# for 2022-Q3, shifting the existing rows down.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q3
#    and shift the existing rows down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()

# Bring over the formatting of the (now shifted) row 3 into the new row 2
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 28
$summary.Cells.Item(2,4).Value = 2.96

# ---------------------------------------------------------------------
# 2. Create the new "2022-Q3" worksheet right after "总计" and fill it
#    with the fund-holding detail, using the existing "2022-Q2" sheet
#    as a formatting template.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q3"

# Copy header row + first 26 data rows (formats + values) from the
# template sheet - gives us correct styles (bold/border/center header,
# bold/border/center column-A index cells) for free.
$template.Range("A1:H27").Copy($newSheet.Range("A1:H27"))
# We need 28 data rows (2022-Q3 has more holdings); clone the last
# template row's formatting twice more for rows 28 and 29.
$template.Range("A27:H27").Copy($newSheet.Range("A28:H28"))
$template.Range("A27:H27").Copy($newSheet.Range("A29:H29"))

# Columns B:G hold text-like values (fund codes with leading zeros,
# names, and numbers formatted as text such as "42.40"); force text
# number format so assigning the values below does not get silently
# re-interpreted as numbers.
$newSheet.Range("B2:G29").NumberFormat = "@"
# ... except G28/G29, which are genuinely numeric zeros in the source.
$newSheet.Range("G28:G29").NumberFormat = "General"

$q3Data = @(
    @("0","163407","兴全沪深300指数增强（LOF）A","42.40","94.96","2.75","1.1660","8"),
    @("1","501202","华泰紫金科技创新3年封闭运作灵活配置混合C","9.12","71.92","3.56","0.3247","6"),
    @("2","159997","天弘中证电子ETF","10.65","99.85","2.19","0.2332","10"),
    @("3","515750","富国中证科技50策略ETF","5.40","99.62","2.66","0.1436","9"),
    @("4","515260","华宝中证电子50ETF","4.62","98.62","2.83","0.1307","10"),
    @("5","159786","银华中证虚拟现实主题ETF","2.41","97.76","5.31","0.1280","5"),
    @("6","515150","富国中证国企一带一路ETF","5.83","99.44","1.99","0.1160","8"),
    @("7","014189","南方专精特新混合A","2.50","83.90","4.21","0.1052","5"),
    @("8","159916","建信深证基本面60ETF","3.55","98.81","2.83","0.1005","9"),
    @("9","159910","嘉实深证基本面120ETF","3.42","99.58","2.27","0.0776","9"),
    @("10","515110","易方达中证国企一带一路ETF","3.60","99.20","1.98","0.0713","8"),
    @("11","515320","华安中证电子50ETF","1.97","98.40","2.82","0.0556","10"),
    @("12","163116","申万菱信中证申万电子行业投资指数（LOF）A","2.02","93.70","2.45","0.0495","9"),
    @("13","159732","华夏国证消费电子主题ETF","1.11","99.40","4.28","0.0475","6"),
    @("14","007230","兴全沪深300指数增强（LOF）C","1.48","94.96","2.75","0.0407","8"),
    @("15","014190","南方专精特新混合C","0.85","83.90","4.21","0.0358","5"),
    @("16","159709","工银瑞信深证物联网50ETF","1.02","98.22","3.27","0.0334","10"),
    @("17","009663","华泰紫金科技创新3年封闭运作灵活配置混合A","0.71","71.92","3.56","0.0253","6"),
    @("18","970043","东吴裕盈一年持有期灵活配置混合A","0.96","52.43","2.54","0.0244","6"),
    @("19","515990","汇添富中证国企一带一路ETF","0.93","98.74","1.98","0.0184","8"),
    @("20","970045","东吴裕盈一年持有期灵活配置混合C","0.44","52.43","2.54","0.0112","6"),
    @("21","970044","东吴裕盈一年持有期灵活配置混合B","0.27","52.43","2.54","0.0069","6"),
    @("22","006906","创金合信鑫收益灵活配置混合E","0.19","52.03","2.60","0.0049","6"),
    @("23","010531","申万菱信中证申万电子行业投资指数（LOF）C","0.19","93.70","2.45","0.0047","9"),
    @("24","006718","国融融盛龙头严选混合A","0.08","60.25","4.55","0.0036","1"),
    @("25","006719","国融融盛龙头严选混合C","0.02","60.25","4.55","0.0009","1"),
    @("26","003749","创金合信鑫收益灵活配置混合A","0.00","52.03","2.60","0","6"),
    @("27","003750","创金合信鑫收益灵活配置混合C","0.00","52.03","2.60","0","6"),
)

for ($i = 0; $i -lt $q3Data.Count; $i++) {
    $row = $i + 2
    $rec = $q3Data[$i]

    $newSheet.Cells.Item($row,1).Value = [int]$rec[0]
    $newSheet.Cells.Item($row,2).Value = $rec[1]
    $newSheet.Cells.Item($row,3).Value = $rec[2]
    $newSheet.Cells.Item($row,4).Value = $rec[3]
    $newSheet.Cells.Item($row,5).Value = $rec[4]
    $newSheet.Cells.Item($row,6).Value = $rec[5]

    if ($row -eq 28 -or $row -eq 29) {
        $newSheet.Cells.Item($row,7).Value = 0
    } else {
        $newSheet.Cells.Item($row,7).Value = $rec[6]
    }

    $newSheet.Cells.Item($row,8).Value = [int]$rec[7]
}

Write-Host "edit complete"
